$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (column D) and Volume(1h) (column E) values.
# D-column prices are plain text (e.g. "61.542.17" / "0.999") in the source sheet,
# so force Text number-format before writing to stop Excel from re-typing them as
# numbers, then restore the default "Normal" style so no stray style id is left on
# the cell (matches the original formatting exactly).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.542.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.377.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.381.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.594.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.512.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("E28").Value = "  +7.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.163"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("E36").Value = "  -6.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.773"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.07%  "
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.354.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  +0.89%  "
